$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Header cell H1 - "Save", using same style as other header cells (e.g. G1)
$ws.Range("H1").Value = "Save"
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122) # xlPasteFormats
$excel.CutCopyMode = $false

# Save column values for rows 2-15
$saveValues = @(0, 1, 0, 0, 0, 0, 1, 1, 1, 1, 0, 0, 1, 0)

for ($i = 0; $i -lt $saveValues.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 8).Value = $saveValues[$i]
}
